$wb = $excel.ActiveWorkbook

# "meta" sheet: add a new key/value row "style" / "default" (row 4), and
# push the previously-blank styled cell down to row 5.
$metaSheet = $wb.Worksheets.Item("meta")

# Copy the formatting of the existing blank key cell (A4, style index 1)
# down to the new blank row (A5) before overwriting A4 with real content.
$metaSheet.Range("A4").Copy()
$metaSheet.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

$metaSheet.Range("A4").Value = "style"
$metaSheet.Range("B4").Value = "default"

$excel.CutCopyMode = 0
